# Updates the crypto price/volume table with the latest scraped values.
# D column = Price text, E column = Volume(1h) text (kept as plain text,
# exactly like the source data - note some "prices" look numeric, e.g.
# "6.19", "0.0419", so we briefly mark the cell as Text (@) before
# assigning, then clear the format again so we don't leave a lasting
# number-format override on the cell).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value2 = $val
    $c.ClearFormats()
}

$ws.Range("D2").Value2 = "64.912.78"
$ws.Range("E2").Value2 = "  +2.07%  "

$ws.Range("D3").Value2 = "3.164.73"
$ws.Range("E3").Value2 = "  +4.18%  "

Set-TextValue "D5" "579.03"
$ws.Range("E5").Value2 = "  +4.43%  "

Set-TextValue "D6" "150.55"
$ws.Range("E6").Value2 = "  +6.97%  "

$ws.Range("E7").Value2 = "  +0.02%  "

$ws.Range("D8").Value2 = "3.164.98"
$ws.Range("E8").Value2 = "  +4.25%  "

$ws.Range("E9").Value2 = "  +2.31%  "

Set-TextValue "D11" "6.19"
$ws.Range("E11").Value2 = "  +0.88%  "

$ws.Range("E12").Value2 = "  +3.62%  "

$ws.Range("E13").Value2 = "  +19.07%  "

Set-TextValue "D14" "37.47"
$ws.Range("E14").Value2 = "  +6.14%  "

$ws.Range("D15").Value2 = "3.682.99"
$ws.Range("E15").Value2 = "  +4.24%  "

$ws.Range("D16").Value2 = "64.987.57"
$ws.Range("E16").Value2 = "  +2.19%  "

$ws.Range("D17").Value2 = "3.165.13"
$ws.Range("E17").Value2 = "  +4.34%  "

$ws.Range("E19").Value2 = "  +1.62%  "

Set-TextValue "D20" "510.28"
$ws.Range("E20").Value2 = "  +8.16%  "

Set-TextValue "D21" "14.84"
$ws.Range("E21").Value2 = "  +6.02%  "

$ws.Range("E22").Value2 = "  +6.91%  "

Set-TextValue "D23" "15.32"
$ws.Range("E23").Value2 = "  +6.07%  "

$ws.Range("E24").Value2 = "  +3.86%  "

Set-TextValue "D25" "84.98"
$ws.Range("E25").Value2 = "  +3.18%  "

$ws.Range("E26").Value2 = "  +0.06%  "

$ws.Range("E27").Value2 = "  +12.85%  "

$ws.Range("E28").Value2 = "  +5.28%  "

Set-TextValue "D29" "2.18"
$ws.Range("E29").Value2 = "  +8.31%  "

$ws.Range("E31").Value2 = "  +15.11%  "

$ws.Range("E32").Value2 = "  +0.08%  "

Set-TextValue "D33" "1.19"
$ws.Range("E33").Value2 = "  +4.13%  "

Set-TextValue "D34" "6.31"
$ws.Range("E34").Value2 = "  +11.84%  "

$ws.Range("E35").Value2 = "  +6.68%  "

Set-TextValue "D36" "55.73"
$ws.Range("E36").Value2 = "  +1.63%  "

$ws.Range("E37").Value2 = "  +10.62%  "

Set-TextValue "D38" "471.79"
$ws.Range("E38").Value2 = "  +7.67%  "

Set-TextValue "D39" "3.12"
$ws.Range("E39").Value2 = "  +14.20%  "

Set-TextValue "D40" "0.0419"
$ws.Range("E40").Value2 = "  +3.59%  "

$ws.Range("E41").Value2 = "  +4.78%  "

$ws.Range("D42").Value2 = "3.061.84"
$ws.Range("E42").Value2 = "  +2.25%  "

$ws.Range("E43").Value2 = "  +1.59%  "

Set-TextValue "D44" "0.282"
$ws.Range("E44").Value2 = "  +5.70%  "

$ws.Range("E45").Value2 = "  +8.83%  "

Set-TextValue "D46" "29.30"
$ws.Range("E46").Value2 = "  +6.27%  "

$ws.Range("D47").Value2 = "0.0₃0604"
$ws.Range("E47").Value2 = "  +19.15%  "

$ws.Range("E49").Value2 = "  +1.62%  "

$ws.Range("E50").Value2 = "  +8.75%  "

Set-TextValue "D51" "119.69"
$ws.Range("E51").Value2 = "  +1.79%  "
